$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pure value corrections in the 07/2025 block (before the row insertion
# shifts anything, so these row numbers are the original ones) ---
$ws.Range("B5").Value = 27892.39    # Dia 4 (07/2025): 26369.89 -> 27892.39
$ws.Range("B21").Value = 12811.82   # Dia 28 (07/2025): 9500.17 -> 12811.82
$ws.Range("B22").Value = 24722.29   # Dia 29 (07/2025): 23549.84 -> 24722.29

# --- Insert the missing "Dia 30" entry for 07/2025, right after Dia 29 (row 22).
# This pushes every following row (old row 23 onward) down by one, which is
# exactly the shift seen for the rest of the sheet in the diff. ---
$ws.Rows(23).Insert()

$ws.Range("A23").Value = 30
$ws.Range("B23").Value = 35906.71
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 2025
$ws.Range("E23").Value = "07/2025"
